# Edit script: add "Player Info" sheet before "ODI Batting", and rename
# MATCH_CARD_LINK columns to MATCH_CODE with just the numeric match code
# stored as text (instead of the full URL).

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "Player Info" sheet as the first sheet ---
$battingSheet = $wb.Worksheets.Item("ODI Batting")

$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"
$playerInfo.Move($battingSheet)
# NOTE: Move() invalidates/repositions existing sheet handles (they track
# position, not identity), so any sheet reference used before this point
# must be re-fetched by name afterwards.

# Header row (bold, bordered, centered - matching the style used by the
# other sheets' header rows)
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$hdr = $playerInfo.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# Data row. The ID looks numeric but must be kept as text, matching how
# the rest of the workbook stores every value (including numeric-looking
# ones) as text. Forcing NumberFormat="@" keeps the value textual but
# also stamps a style onto the cell, so we strip that back off with
# ClearFormats() once the text value has been recorded.
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4390"
$playerInfo.Range("A2").ClearFormats()
$playerInfo.Range("B2").Value = "Jayant Yadav"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# --- 2. Update "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$ws1 = $wb.Worksheets.Item("ODI Batting")
$ws1.Range("D1").Value = "MATCH_CODE"

$ws1.Range("D2").NumberFormat = "@"
$ws1.Range("D2").Value = "3955"
$ws1.Range("D2").ClearFormats()

$ws1.Range("D3").NumberFormat = "@"
$ws1.Range("D3").Value = "4529"
$ws1.Range("D3").ClearFormats()

# --- 3. Update "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$ws2 = $wb.Worksheets.Item("ODI Bowling")
$ws2.Range("B1").Value = "MATCH_CODE"

$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "3955"
$ws2.Range("B2").ClearFormats()

$ws2.Range("B3").NumberFormat = "@"
$ws2.Range("B3").Value = "4529"
$ws2.Range("B3").ClearFormats()
